# Generate Report for Handback
# The handback report re-sorts the two tracked files (560101b1... and
# f30d5173...) on every sheet: the f30d5173 item is now "Handed back: in
# sync with en-US" (and has gained its handback target/file columns), and
# moves into the row that 560101b1 used to occupy; 560101b1 drops into the
# row f30d5173 used to occupy, keeping its "Ready for handoff" status.

$wb = $excel.ActiveWorkbook

$mdUrl560101 = "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/560101b1-a1ef-4878-b250-f85ebe891b31.md"
$mdUrlF30d = "https://github.com/OpenLocalizationTest/oltest/blob/a456910e3c573c22ae05a56b82c899301028891f/e2e/f30d5173-3193-4915-b48c-f81210d73ab0.md"

$zhXlfUrl560101 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c162d1a28e6467a7c4e80d32583a4e5e14fc0adb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/560101b1-a1ef-4878-b250-f85ebe891b31.843ee1106a8550accee87f7b8a8c33b31aba932a.zh-cn.xlf"
$zhXlfUrlF30d = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d603c81b1693c7eaf6df65be6220e089f5f558d5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f30d5173-3193-4915-b48c-f81210d73ab0.23c9c913752ae51596a2004e31e078ee7e3796c1.zh-cn.xlf"

$deXlfUrl560101 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9abb8d041e474ecd509046179e256280413fabe3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/560101b1-a1ef-4878-b250-f85ebe891b31.843ee1106a8550accee87f7b8a8c33b31aba932a.de-de.xlf"
$deXlfUrlF30d = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e0bf6a265f1635d7e75c25bd9225ccd1fb3a8dbf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f30d5173-3193-4915-b48c-f81210d73ab0.23c9c913752ae51596a2004e31e078ee7e3796c1.de-de.xlf"

$md560101 = "560101b1-a1ef-4878-b250-f85ebe891b31.md"
$mdF30d = "f30d5173-3193-4915-b48c-f81210d73ab0.md"

$zhXlf560101 = "560101b1-a1ef-4878-b250-f85ebe891b31.843ee1106a8550accee87f7b8a8c33b31aba932a.zh-cn.xlf"
$zhXlfF30d = "f30d5173-3193-4915-b48c-f81210d73ab0.23c9c913752ae51596a2004e31e078ee7e3796c1.zh-cn.xlf"

$deXlf560101 = "560101b1-a1ef-4878-b250-f85ebe891b31.843ee1106a8550accee87f7b8a8c33b31aba932a.de-de.xlf"
$deXlfF30d = "f30d5173-3193-4915-b48c-f81210d73ab0.23c9c913752ae51596a2004e31e078ee7e3796c1.de-de.xlf"

$handedBack = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"
$include = "Include"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsOverview.Range("D2").Value = "2016-34-20 16:34:11"

$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff
$wsOverview.Range("D3").Value = "2016-33-20 16:33:49"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl560101, "", "", $mdF30d)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrlF30d, "", "", $md560101)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("E2").Value = "2016-03-20 16:34:08"
$wsZh.Range("F2").Value = $mdF30d
$wsZh.Range("G2").Value = $zhXlfF30d
$wsZh.Range("H2").Value = "2016-03-20 16:34:27"
$wsZh.Range("I2").Value = $include

$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("E3").Value = "2016-03-20 16:33:46"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = $include

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrlF30d, "", "", $mdF30d)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdUrlF30d, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrlF30d, "", "", $zhXlfF30d)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrlF30d, "", "", $mdF30d)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrlF30d, "", "", $zhXlfF30d)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl560101, "", "", $md560101)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl560101, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl560101, "", "", $zhXlf560101)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("E2").Value = "2016-03-20 16:34:11"
$wsDe.Range("F2").Value = $mdF30d
$wsDe.Range("G2").Value = $deXlfF30d
$wsDe.Range("H2").Value = "2016-03-20 16:34:32"
$wsDe.Range("I2").Value = $include

$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("E3").Value = "2016-03-20 16:33:49"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = $include

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrlF30d, "", "", $mdF30d)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdUrlF30d, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrlF30d, "", "", $deXlfF30d)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrlF30d, "", "", $mdF30d)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrlF30d, "", "", $deXlfF30d)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl560101, "", "", $md560101)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl560101, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl560101, "", "", $deXlf560101)
